# Apply row-window corrections to CryCompanywiseStockReport data.
# Upstream sent misaligned columns B/C/E/F/G for the rows below; each group
# of consecutive item rows needs those columns rotated up by one row (the
# last row in a group receives the first row's original values) while the
# serial number (A) and base price (D) columns stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Rotate-RowWindow {
    param($Worksheet, [int[]]$Rows)

    $cols = @("B", "C", "E", "F", "G")

    # Snapshot the current values for every row in the window first so the
    # writes below never read back an already-overwritten cell.
    $snapshot = @{}
    foreach ($r in $Rows) {
        $rowValues = @{}
        foreach ($c in $cols) {
            $rowValues[$c] = $Worksheet.Range("$c$r").Value()
        }
        $snapshot[$r] = $rowValues
    }

    $count = $Rows.Count
    for ($i = 0; $i -lt $count; $i++) {
        $targetRow = $Rows[$i]
        $sourceRow = $Rows[($i + 1) % $count]
        $sourceValues = $snapshot[$sourceRow]
        foreach ($c in $cols) {
            $Worksheet.Range("$c$targetRow").Value = $sourceValues[$c]
        }
    }
}

Rotate-RowWindow $ws @(306, 307)
Rotate-RowWindow $ws @(339, 340)
Rotate-RowWindow $ws @(343, 344, 345)
Rotate-RowWindow $ws @(348, 349)
Rotate-RowWindow $ws @(365, 366)
Rotate-RowWindow $ws @(375, 376)
Rotate-RowWindow $ws @(382, 383)
Rotate-RowWindow $ws @(393, 394)
Rotate-RowWindow $ws @(412, 413)
Rotate-RowWindow $ws @(424, 425)
Rotate-RowWindow $ws @(572, 573)
Rotate-RowWindow $ws @(579, 580)
Rotate-RowWindow $ws @(583, 584)
Rotate-RowWindow $ws @(586, 587)
Rotate-RowWindow $ws @(680, 681)
Rotate-RowWindow $ws @(702, 703)
Rotate-RowWindow $ws @(713, 714)
Rotate-RowWindow $ws @(865, 866)

Write-Output "Rotated columns B,C,E,F,G for row groups: 306,307; 339,340; 343,344,345; 348,349; 365,366; 375,376; 382,383; 393,394; 412,413; 424,425; 572,573; 579,580; 583,584; 586,587; 680,681; 702,703; 713,714; 865,866"
